$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerate column G ("K") values for rows 2-16 (Strike# -> K)
$kValues = @{
    2  = 5
    3  = 3
    4  = 4
    5  = 4
    6  = 4
    7  = 7
    8  = 5
    9  = 7
    10 = 7
    11 = 6
    12 = 8
    13 = 2
    14 = 4
    15 = 9
    16 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
